$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '66.081.45'
$ws.Range("E2").Value = '  +7.28%  '
$ws.Range("D3").Value = '3.016.13'
$ws.Range("E3").Value = '  +4.22%  '
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").Value = '586.32'
$ws.Range("E5").Value = '  +3.03%  '
$ws.Range("D6").Value = '154.62'
$ws.Range("E6").Value = '  +7.22%  '
$ws.Range("E7").Value = '  -0.19%  '
$ws.Range("D8").Value = '3.011.75'
$ws.Range("E8").Value = '  +4.06%  '
$ws.Range("E9").Value = '  +2.51%  '
$ws.Range("E10").Value = '  -0.06%  '
$ws.Range("D11").Value = '0.154'
$ws.Range("E11").Value = '  +4.87%  '
$ws.Range("D12").Value = '0.451'
$ws.Range("E12").Value = '  +4.56%  '
$ws.Range("E13").Value = '  +3.48%  '
$ws.Range("E14").Value = '  +7.03%  '
$ws.Range("E15").Value = '  +0.75%  '
$ws.Range("D16").Value = '66.040.38'
$ws.Range("E16").Value = '  +7.20%  '
$ws.Range("D17").Value = '3.510.93'
$ws.Range("E17").Value = '  +4.08%  '
$ws.Range("D18").Value = '6.98'
$ws.Range("E18").Value = '  +6.29%  '
$ws.Range("D19").Value = '3.011.55'
$ws.Range("E19").Value = '  +4.11%  '
$ws.Range("D20").Value = '457.70'
$ws.Range("E20").Value = '  +5.63%  '
$ws.Range("D21").Value = '13.84'
$ws.Range("E21").Value = '  +5.53%  '
$ws.Range("D22").Value = '0.686'
$ws.Range("E22").Value = '  +4.31%  '
$ws.Range("E23").Value = '  +7.69%  '
$ws.Range("D24").Value = '81.82'
$ws.Range("E24").Value = '  +3.10%  '
$ws.Range("D25").Value = '12.65'
$ws.Range("E25").Value = '  +5.19%  '
$ws.Range("E26").Value = '  +11.78%  '
$ws.Range("D27").Value = '10.73'
$ws.Range("E27").Value = '  +7.27%  '
$ws.Range("E28").Value = '  -0.04%  '
$ws.Range("D29").Value = '2.44'
$ws.Range("E29").Value = '  +19.04%  '
$ws.Range("D30").Value = '7.89'
$ws.Range("E30").Value = '  +12.10%  '
$ws.Range("E31").Value = '  +4.02%  '
$ws.Range("E32").Value = '  -2.50%  '
$ws.Range("E33").Value = '  +4.93%  '
$ws.Range("D34").Value = '27.05'
$ws.Range("E34").Value = '  +6.02%  '
$ws.Range("D35").Value = '1.00'
$ws.Range("E35").Value = '  +0.11%  '
$ws.Range("E36").Value = '  +3.28%  '
$ws.Range("D37").Value = '5.79'
$ws.Range("E37").Value = '  +7.27%  '
$ws.Range("D38").Value = '2.15'
$ws.Range("E38").Value = '  +11.22%  '
$ws.Range("D39").Value = '45.59'
$ws.Range("E39").Value = '  +15.20%  '
$ws.Range("D40").Value = '49.40'
$ws.Range("E40").Value = '  +1.03%  '
$ws.Range("D41").Value = '2.95'
$ws.Range("E41").Value = '  +4.20%  '
$ws.Range("E42").Value = '  +6.31%  '
$ws.Range("E43").Value = '  +13.32%  '
$ws.Range("D44").Value = '8.47'
$ws.Range("E44").Value = '  +2.91%  '
$ws.Range("D45").Value = '387.49'
$ws.Range("E45").Value = '  +11.49%  '
$ws.Range("D46").Value = '2.782.74'
$ws.Range("E46").Value = '  +2.83%  '
$ws.Range("E47").Value = '  +5.38%  '
$ws.Range("D48").Value = '134.79'
$ws.Range("E48").Value = '  +1.73%  '
$ws.Range("D50").Value = '23.49'
$ws.Range("E50").Value = '  +8.54%  '
$ws.Range("E51").Value = '  +3.20%  '
